$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string "FAPs" used by new rows (will be appended to shared strings table)

$data = @(
    @(2, "ECs", "Tnfsf13", "Tnfrsf13b", "M1", 2, 1, 0.7536765000000001, 1.507353, 0.07241044948256109, 0.05015328104483915, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 8.463944864614001, 50.783669187684, 0.02977078241293467, 0.02061998548483333),
    @(3, "ECs", "Tnfsf13", "Tnfrsf13b", "M2", 2, 1, 0.7536765000000001, 1.507353, 0.07241044948256109, 0.05015328104483915, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 12.1226169375385, 72.73570162523102, 0.04263966706962641, 0.02953329556000581),
    @(4, "M1", "Tnfsf13", "Tnfrsf13b", "M1", 2, 0.6666666666666666, 0.1204533333333333, 0.36136, 0.01157271058383726, 0.01202332143722344, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 1.352716148231111, 12.17444533408, 0.00475799627237861, 0.004943260108812847),
    @(5, "M1", "Tnfsf13", "Tnfrsf13b", "M2", 2, 0.6666666666666666, 0.1204533333333333, 0.36136, 0.01157271058383726, 0.01202332143722344, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 1.937448784524445, 17.43703906072, 0.006814714311458651, 0.007080061328410596),
    @(6, "M2", "Tnfsf13", "Tnfrsf13b", "M1", 3, 1, 3.418481, 10.255443, 0.3284350059443208, 0.341223399574173, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 38.39025723202267, 345.512315088204, 0.135032542521561, 0.1402903538856098),
    @(7, "M2", "Tnfsf13", "Tnfrsf13b", "M2", 3, 1, 3.418481, 10.255443, 0.3284350059443208, 0.341223399574173, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 54.98504420829568, 494.8653978746611, 0.1934024634227597, 0.2009330456885631),
    @(8, "Neutro", "Tnfsf13", "Tnfrsf13b", "M1", 3, 1, 4.597592, 13.792776, 0.44171962805982, 0.4589190263438705, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 51.63192058925866, 464.6872853033279, 0.181608304166906, 0.1886796529516029),
    @(9, "Neutro", "Tnfsf13", "Tnfrsf13b", "M2", 3, 1, 4.597592, 13.792776, 0.44171962805982, 0.4589190263438705, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 73.95062291459467, 665.555606231352, 0.260111323892914, 0.2702393733922675),
    @(10, "sCs", "Tnfsf13", "Tnfrsf13b", "M1", 3, 1, 1.101608333333333, 3.304825, 0.1058384526655689, 0.1099595231037524, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 12.37129218667778, 111.3416296801, 0.04351434865747077, 0.045208682723897),
    @(11, "sCs", "Tnfsf13", "Tnfrsf13b", "M2", 3, 1, 1.101608333333333, 3.304825, 0.1058384526655689, 0.1099595231037524, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 17.71897603308611, 159.470784297775, 0.06232410400809811, 0.0647508403798554),
    @(12, "FAPs", "Tnfsf13", "Tnfrsf13b", "M1", 2, 1, 0.416583, 0.8331660000000001, 0.04002375326389206, 0.02772144849614155, 3, 1, 11.23020933333333, 33.690628, 0.4111393124289954, 0.4111393124289954, 4.678314294708001, 28.069885768248, 0.01645533839774434, 0.01139737727423944),
    @(13, "FAPs", "Tnfsf13", "Tnfrsf13b", "M2", 2, 1, 0.416583, 0.8331660000000001, 0.04002375326389206, 0.02772144849614155, 3, 1, 16.08464233333333, 48.253927, 0.5888606875710045, 0.5888606875710045, 6.700588557147001, 40.20353134288201, 0.02356841486614771, 0.0163240712219021),
)

foreach ($row in $data) {
    $r = [int]$row[0]
    for ($col = 1; $col -le 20; $col++) {
        $val = $row[$col]
        if ($col -le 4) {
            $ws.Cells.Item($r, $col).Value = $val
        } else {
            $ws.Cells.Item($r, $col).Value = [double]$val
        }
    }
}

Write-Host "done"
